# Update price return for index comparison
# Applies updated Market Price / Market Value / Capital Gains / Total Return
# figures (and the downstream summary + index-comparison sheets) to reflect
# a revised market price for AY and ENB.TO.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "Stock log"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Stock log")

# Row 2 (AY buy lot #1)
$ws1.Range("G2").Value = 21.59
$ws1.Range("K2").Value = 14832.33
$ws1.Range("L2").Value = -8.32
$ws1.Range("N2").Value = -4.72

# Row 3 (AY buy lot #2)
$ws1.Range("G3").Value = 21.59
$ws1.Range("K3").Value = 8636
$ws1.Range("L3").Value = -12.45
$ws1.Range("N3").Value = -9.17

# Row 5 (AY sell lot)
$ws1.Range("G5").Value = 21.59
$ws1.Range("K5").Value = -2159
$ws1.Range("L5").Value = -6.13
$ws1.Range("N5").Value = -4.3

# Row 6 (ENB.TO buy lot)
$ws1.Range("G6").Value = 35.73
$ws1.Range("K6").Value = 3573
$ws1.Range("L6").Value = 5.87
$ws1.Range("N6").Value = 10.23

# ---------------------------------------------------------------------
# Sheet 2: "Portfolio Summary"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Portfolio Summary")

# Row 2 (AY)
$ws2.Range("D2").Value = 21309.33
$ws2.Range("G2").Value = -10.27
$ws2.Range("H2").Value = -6.62

# Row 4 (ENB.TO)
$ws2.Range("D4").Value = 3573
$ws2.Range("G4").Value = 5.87
$ws2.Range("H4").Value = 10.21

# ---------------------------------------------------------------------
# Sheet 3: "Total Return"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Total Return")

# Row 2 (My Portfolio)
$ws3.Range("C2").Value = -1.57
$ws3.Range("D2").Value = 0.84

# Row 3 (Global X Super Dividend ETF)
$ws3.Range("C3").Value = -1.99

# Row 4 (S&P 500)
$ws3.Range("C4").Value = 25.85
